# Move the "simulation_times" sheet's single row of values into a new
# row 17 ("simtime") on the "optimization_parameters" sheet, then remove
# the now-redundant "simulation_times" sheet, leaving "network_b" as the
# last tab. Finally, leave "optimization_parameters" as the active sheet
# with A17 selected (mirrors the author's post-edit selection state).

$wb = $excel.ActiveWorkbook

$wsTimes = $wb.Worksheets.Item("simulation_times")
$wsOpt   = $wb.Worksheets.Item("optimization_parameters")

# Copy the time-series values (B1:V1 on simulation_times) down to row 17
# of optimization_parameters, labeled "simtime" in column A.
$wsOpt.Cells.Item(17, 1).Value2 = "simtime"
for ($col = 2; $col -le 22; $col++) {
    $wsOpt.Cells.Item(17, $col).Value2 = $wsTimes.Cells.Item(1, $col).Value2
}

# The simulation_times sheet is no longer needed now that its data lives
# on optimization_parameters.
[void]$wsTimes.Delete()

# Make optimization_parameters the active tab with A17 selected.
[void]$wsOpt.Activate()
[void]$wsOpt.Range("A17").Select()
